$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New permission rows to support importing multiple "work basis" records
# (sift enforcement staff / inspected org, desktop-client import/export flow).
$rows = @(
    @{ Id = 320; Name = "跳转到筛选执法人员和被检查机构页面"; Method = "toSiftPage" },
    @{ Id = 321; Name = "筛选其他组成员"; Method = "siftAePeoples" },
    @{ Id = 322; Name = "筛选被检查机构"; Method = "siftAeedOrg" },
    @{ Id = 323; Name = "保存筛选结果"; Method = "toAdminEnforceInitPageFromSift" },
    @{ Id = 324; Name = "筛选结果列表页面"; Method = "toSiftResultList" },
    @{ Id = 325; Name = "从筛选结果跳转行政执法登记页面"; Method = "toAdminEnforceInitPageFromSiftResult" },
    @{ Id = 326; Name = "导出单机版初始化文件"; Method = "generateDesktopClientInitialFile" },
    @{ Id = 327; Name = "跳转行政执法的工作检查记录单机版数据导入界面"; Method = "toImportDesktopClientFinalData" },
    @{ Id = 328; Name = "导入行政执法的工作检查记录单机版数据"; Method = "importDesktopClientFinalData" }
)

$startRow = 118
$r = $startRow
foreach ($row in $rows) {
    $ws.Range("A$r").Value = $row.Id
    $ws.Range("B$r").Value = "net.sf.jguard.core.authorization.permissions.URLPermission"
    $ws.Range("C$r").Value = $row.Name
    $ws.Range("D$r").Value = $row.Method
    $ws.Range("E$r").Value = 5
    $ws.Range("F$r").Value = "/AdminEnforceManagerAction.do"
    $ws.Range("G$r").Formula = "=""INSERT INTO JG_PERMISSION (ID, CLASS, NAME, ACTIONS, DOMAIN_ID) VALUES (""&A$r&"", '""&B$r&""','""&C$r&""','""&F$r&""?method=""&D$r&""&*,ANY',""&E$r&"");"""
    $r = $r + 1
}

$ws.Range("C118").Font.Size = 11

$ws.Range("A118").Select()
